# Add two new roster entries ("fails4" and "fails3ab") to the "blackboard"
# sheet, mirroring the pattern of the existing rows (2-7), and move the
# active selection to D14 as in the post-edit workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("blackboard")

$dateValue = 44044.922222222223
$dateFormat = "m/d/yy h:mm"

# Row 8: fails4
$ws.Cells.Item(8, 1).Value = "Name"
$ws.Cells.Item(8, 2).Value = "Name"
$ws.Cells.Item(8, 3).Value = "fails4"
$ws.Cells.Item(8, 4).Value = "fails4"
$ws.Cells.Item(8, 5).Value = $dateValue
$ws.Cells.Item(8, 5).NumberFormat = $dateFormat
$ws.Cells.Item(8, 6).Value = "Yes"

# Row 9: fails3ab
$ws.Cells.Item(9, 1).Value = "Name"
$ws.Cells.Item(9, 2).Value = "Name"
$ws.Cells.Item(9, 3).Value = "fails3ab"
$ws.Cells.Item(9, 4).Value = "fails3ab"
$ws.Cells.Item(9, 5).Value = $dateValue
$ws.Cells.Item(9, 5).NumberFormat = $dateFormat
$ws.Cells.Item(9, 6).Value = "Yes"

# Move the selection/active cell as recorded after the edit.
$null = $ws.Range("D14").Select()
